# ppPlaceholderDate, per the PowerPoint object model's PpPlaceholderType enum.
$ppPlaceholderDate = 16

$p = $ppt.ActivePresentation

# Update the "datetimeFigureOut" date placeholder text on the slide master
# and on every slide layout from "5/7/18" to "5/13/18".
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
        $shp.TextFrame.TextRange.Text = "5/13/18"
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $shp.TextFrame.TextRange.Text = "5/13/18"
        }
    }
}

# Remove the "Straight Connector 33" connector shape from slide 1.
$s = $p.Slides.Item(1)
$conn = $s.Shapes.Item("Straight Connector 33")
$conn.Delete()
